$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "a"
$ws.Range("A6").Value = "aa"
$ws.Range("C6").Value = "a"
$ws.Range("B7").Value = "a"
$ws.Range("C8").Value = "a"

$ws.Range("A6").Select()
